$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Facility utilisation (E) and Fees paid to UKEF currency (F) values for rows 2-6
$ws.Range("E2").Value = 800000
$ws.Range("E3").Value = 800000
$ws.Range("E4").Value = 800000
$ws.Range("E5").Value = 800000
$ws.Range("E6").Value = 800000

$ws.Range("F2").Value = 761579.37
$ws.Range("F3").Value = 761579.37
$ws.Range("F4").Value = 761579.37

# Swap values of G5 and H5
$ws.Range("G5").Value = 456
$ws.Range("H5").Value = 3938753.8

# Update the active selection on the sheet
$ws.Range("E2:H6").Select()
